# Update cryptocurrency Price (D) and Volume(1h) (E) columns per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.930.61'
$ws.Range('E2').Value = '  -0.08%  '

$ws.Range('D3').Value = '3.382.05'
$ws.Range('E3').Value = '  -0.39%  '

$ws.Range('D5').Value = "'570.79"
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.38%  '

$ws.Range('D6').Value = "'141.69"
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.81%  '

$ws.Range('E7').Value = '  +0.00%  '

$ws.Range('E8').Value = '  -0.02%  '

$ws.Range('E9').Value = '  +0.62%  '

$ws.Range('E10').Value = '  -1.56%  '

$ws.Range('D11').Value = "'0.387"
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.87%  '

$ws.Range('D12').Value = '3.961.38'
$ws.Range('E12').Value = '  -0.28%  '

$ws.Range('E13').Value = '  +1.83%  '

$ws.Range('D14').Value = "'27.78"
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.83%  '

$ws.Range('E15').Value = '  -0.09%  '

$ws.Range('D16').Value = '3.390.23'
$ws.Range('E16').Value = '  -0.04%  '

$ws.Range('D17').Value = '61.044.93'
$ws.Range('E17').Value = '  +0.00%  '

$ws.Range('E18').Value = '  -3.53%  '

$ws.Range('D19').Value = "'13.62"
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -3.86%  '

$ws.Range('D20').Value = "'8.95"
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -1.75%  '

$ws.Range('D21').Value = "'382.83"
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.83%  '

$ws.Range('D22').Value = "'75.10"
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +2.64%  '

$ws.Range('E23').Value = '  -2.43%  '

$ws.Range('E24').Value = '  +0.09%  '

$ws.Range('E25').Value = '  -2.75%  '

$ws.Range('D26').Value = '3.521.55'
$ws.Range('E26').Value = '  -0.31%  '

$ws.Range('E27').Value = '  +1.17%  '

$ws.Range('D28').Value = "'0.998"
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.26%  '

$ws.Range('D29').Value = "'7.26"
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -1.62%  '

$ws.Range('D30').Value = "'7.96"
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -2.49%  '

$ws.Range('E31').Value = '  -0.56%  '

$ws.Range('E32').Value = '  -0.04%  '

$ws.Range('E33').Value = '  -4.26%  '

$ws.Range('D34').Value = "'23.22"
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -2.87%  '

$ws.Range('D35').Value = "'6.94"
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.65%  '

$ws.Range('D36').Value = "'166.45"
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.29%  '

$ws.Range('D37').Value = '3.415.09'
$ws.Range('E37').Value = '  -0.16%  '

$ws.Range('E38').Value = '  -2.08%  '

$ws.Range('E39').Value = '  -3.88%  '

$ws.Range('D40').Value = "'0.0770"
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -1.73%  '

$ws.Range('D41').Value = "'26.81"
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -1.15%  '

$ws.Range('E42').Value = '  +0.05%  '

$ws.Range('E43').Value = '  -1.21%  '

$ws.Range('E44').Value = '  -2.66%  '

$ws.Range('E45').Value = '  -2.35%  '

$ws.Range('E46').Value = '  -0.72%  '

$ws.Range('D47').Value = '2.449.75'
$ws.Range('E47').Value = '  -3.74%  '

$ws.Range('D48').Value = "'22.92"
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.02%  '

$ws.Range('D49').Value = "'6.71"
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -2.46%  '

$ws.Range('E50').Value = '  +8.87%  '

$ws.Range('E51').Value = '  +1.34%  '
